$d = $word.ActiveDocument

# Locate the "Author Contributions" heading paragraph.
$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Author Contributions*") {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -eq -1) {
    throw "Could not find 'Author Contributions' heading paragraph"
}

# The author-contribution entries run from right after the heading to the
# end of the document body (they are the last paragraphs before sectPr).
$lastIndex = $d.Paragraphs.Count
$startPara = $d.Paragraphs.Item($headingIndex + 1)
$endPara = $d.Paragraphs.Item($lastIndex)

$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()

# Re-insert the paragraphs in the new order, with corrected initials
# (the ambiguous "C.R." split into "C.Ro." for Cristiane Rocha and
# "C.Ru." for Christian Rummel).
$headingPara = $d.Paragraphs.Item($headingIndex)
$insPoint = $d.Range($headingPara.Range.End, $headingPara.Range.End)

$newParagraphsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:i/></w:rPr><w:t>Performed the analysis</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>J.B.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>M.L.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>Y.D.vdW.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Conceived and designed the analysis</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>J.B.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>M.L.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>Y.D.vdW.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Analyzed the data</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>C.Ro.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>C.Y.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>F.C.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>F.P.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>G.S.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>J.B.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>J.D.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>K.L.P.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>K.Z.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>M.L.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>M.R.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>T.M.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>T.P.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>Y.D.vdW.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Cohort co-investigator</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>J.B.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>M.L.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>N.J.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>P.T.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Contributed data or analysis tools</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>J.B.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>M.L.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>Y.D.vdW.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Wrote the paper</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>J.B.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>M.L.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Cohort PI</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>Y.D.vdW.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Collected the data</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>B.G.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>C.M.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>C.Ro.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>C.Ru.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>C.Y.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>F.C.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>F.P.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>G.S.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>I.D.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>J.B.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>J.D.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>K.L.P.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>K.Z.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>L.P.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>M.L.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>M.R.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>O.A.vdH.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>R.M.DB.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>R.W.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>S.a-B.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>T.M.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>T.P.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>Y.D.vdW.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Read, edited and approved the paper</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>B.G.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>C.M.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>C.Ro.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>C.Ru.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>C.Y.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>F.C.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>F.P.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>G.S.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>I.D.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>J.D.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>K.L.P.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>K.Z.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>L.P.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>M.R.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>N.J.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>O.A.vdH.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>P.T.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>R.M.DB.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>R.W.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>S.a-B.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>T.M.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>T.P.</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>Y.D.vdW.</w:t></w:r></w:p>'

$insPoint.InsertXML($newParagraphsXml)

Write-Output "Author Contributions section rebuilt; paragraph count: $($d.Paragraphs.Count)"
